$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 (Segunda / 50 / 1500 / 1500 / 1500 / 1500 / 44466) currently holds the data that,
# per the update, needs to end up duplicated (with edits) across the new rows 38-40 while
# row 37 itself is rewritten with fresh data. Copy the existing row 37 formatting/values
# down into rows 38, 39 and 40 first, before anything in row 37 is touched.
$ws.Range("A37:R37").Copy($ws.Range("A38:R38"))
$ws.Range("A37:R37").Copy($ws.Range("A39:R39"))
$ws.Range("A37:R37").Copy($ws.Range("A40:R40"))

# Row 38: new "Segunda" entry for the 2021-10-22 sample
$ws.Cells.Item(38, 4).Value = 44491
$ws.Cells.Item(38, 9).Value = "Segunda"
$ws.Cells.Item(38, 10).Value = 300
$ws.Cells.Item(38, 11).Value = 1000
$ws.Cells.Item(38, 12).Value = 1000
$ws.Cells.Item(38, 13).Value = 1000
$ws.Cells.Item(38, 16).Value = 1000

# Row 39: keeps the 2021-09-27 date, becomes the "Primera" entry
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 300
$ws.Cells.Item(39, 11).Value = 2000
$ws.Cells.Item(39, 12).Value = 2000
$ws.Cells.Item(39, 13).Value = 2000
$ws.Cells.Item(39, 16).Value = 2000

# Row 40: unchanged copy of the original row 37 data (Segunda / 50 / 1500.../ 44466)

# Row 37 is rewritten with the new "Primera" sample for 2021-10-22
$ws.Cells.Item(37, 4).Value = 44491
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 400
$ws.Cells.Item(37, 11).Value = 1300
$ws.Cells.Item(37, 12).Value = 1300
$ws.Cells.Item(37, 13).Value = 1300
$ws.Cells.Item(37, 16).Value = 1300

# Row 36 is rewritten with the new "Extra" sample for 2021-10-22
$ws.Cells.Item(36, 4).Value = 44491
$ws.Cells.Item(36, 9).Value = "Extra"
$ws.Cells.Item(36, 10).Value = 250
$ws.Cells.Item(36, 11).Value = 1500
$ws.Cells.Item(36, 12).Value = 1500
$ws.Cells.Item(36, 13).Value = 1500
$ws.Cells.Item(36, 16).Value = 1500
